{"js": "// Civil/Plaint.docx edit:\n//  1. Title paragraph \"Plaint\" switches from the Heading1 style to Title.\n//  2. The sample body paragraph is replaced with the full Plaint template\n//     text (a series of paragraphs, some blank, making up the pleading).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Re-style the heading paragraph.\nconst titlePara = paragraphs.items[0];\ntitlePara.style = \"Title\";\n\n// 2) Replace the second paragraph's content with the full template, each\n//    line becoming its own paragraph (joined with carriage returns so Word\n//    splits them into separate paragraphs instead of line breaks).\nconst lines = [\n  \"IN THE COURT OF THE CIVIL JUDGE AT [City]\",\n  \"\",\n  \"PLAINT\",\n  \"\",\n  \"Plaintiff: [Name & Address]\",\n  \"Defendant: [Name & Address]\",\n  \"\",\n  \"The plaintiff respectfully submits:\",\n  \"1. That the plaintiff and defendant entered into [contract/transaction] on [date].\",\n  \"2. That the defendant failed to [action/obligation].\",\n  \"\",\n  \"Prayer:\",\n  \"a) Decree for \\u20b9[Amount]\",\n  \"b) Costs and interest\",\n  \"\",\n  \"Filed on: [Date]\",\n  \"[Signature]\",\n];\n\nconst bodyPara = paragraphs.items[1];\nbodyPara.getRange().insertText(lines.join(\"\\r\"), \"Replace\");\n\nawait context.sync();\n", "ps1": "# Civil/Plaint.docx edit:\n#  1. Title paragraph \"Plaint\" switches from the Heading1 style to Title.\n#  2. The sample body paragraph is replaced with the full Plaint template\n#     text (a series of paragraphs, some blank, making up the pleading).\n\n$d = $word.ActiveDocument\n\n# 1) Re-style the heading paragraph.\n$titlePara = $d.Paragraphs(1)\n$titlePara.Style = \"Title\"\n\n# 2) Replace the second paragraph's content with the full template, each\n#    line becoming its own paragraph (joined with carriage returns so Word\n#    splits them into separate paragraphs instead of line breaks).\n$lines = @(\n    \"IN THE COURT OF THE CIVIL JUDGE AT [City]\",\n    \"\",\n    \"PLAINT\",\n    \"\",\n    \"Plaintiff: [Name & Address]\",\n    \"Defendant: [Name & Address]\",\n    \"\",\n    \"The plaintiff respectfully submits:\",\n    \"1. That the plaintiff and defendant entered into [contract/transaction] on [date].\",\n    \"2. That the defendant failed to [action/obligation].\",\n    \"\",\n    \"Prayer:\",\n    \"a) Decree for \u20b9[Amount]\",\n    \"b) Costs and interest\",\n    \"\",\n    \"Filed on: [Date]\",\n    \"[Signature]\"\n)\n\n$bodyPara = $d.Paragraphs(2)\n$bodyPara.Range.Text = ($lines -join \"`r\")\n"}
